# Weekly driver report update for 2025-04-19
# Refreshes the Bad Drivers and Good Drivers tables on the
# "Driver Summary" sheet with this week's sample/report data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural cleanup -----------------------------------------------
# One fewer "Bad Drivers" entry this week: drop the AX211 23.60.1.2 row
# (row 4); the Totals row shifts up to row 4.
$ws.Rows.Item(4).Delete()

# Three AX211 entries dropped out of the "Good Drivers" table this week;
# after the delete above they sit at row 13 (deleting row 13 three times
# removes all three, since each delete shifts the next one up into it).
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()

# Remaining six "Good Drivers" rows (12-17) are currently in this order:
#   23.100.0.4, 22.80.0.9, 22.50.1.1, 21.110.3.2, 21.70.0.6, 21.60.2.1
# Reorder them to this week's report order:
#   21.60.2.1, 22.50.1.1, 23.100.0.4, 22.80.0.9, 21.110.3.2, 21.70.0.6
$ws.Rows.Item(17).Copy()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(18).Delete()

$ws.Rows.Item(15).Copy()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(16).Delete()

# --- Bad Drivers table ---------------------------------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.0.7"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 97

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 4

# --- Good Drivers table ---------------------------------------------------
# Sample counts refreshed for this week.
$ws.Range("B12").Value = 56018
$ws.Range("B13").Value = 34244
$ws.Range("B14").Value = 442178
$ws.Range("B15").Value = 77849
$ws.Range("B16").Value = 59673
$ws.Range("B17").Value = 113652

# 21.70.0.6's driver vintage is now the date that used to belong to
# 21.60.2.1 ("2019-12-14"); grab it (as text, before it is overwritten)
# via copy/paste so the date-like string isn't re-parsed into a real date.
$ws.Range("E12").Copy()
$ws.Range("E17").PasteSpecial()

# 21.60.2.1 and 22.50.1.1 no longer report a driver vintage this week.
$ws.Range("E12").Value = 0
$ws.Range("E13").Value = 0

# --- Formatting -------------------------------------------------------
# Column A narrows by one character (45 -> 44).
$ws.Columns.Item(1).ColumnWidth = 43.166666666666664
